# Applies the edits described by the commit diff to the PDF manifest workbook:
# fills in a few "Finished"/"Notes" cells for the 2007-2010 rows and scrolls /
# reselects the sheet the way the author left it.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Row 25 (2007, table 3): mark finished and reuse the standard "where do
# other expenses come from?" note.
$ws.Range("H25").Value = "X"
$ws.Range("I25").Value = "Where do other expenses come from?"

# Row 27 (2008, table 3): mark finished and flag the mammals subtotal issue.
$ws.Range("H27").Value = "X"
$ws.Range("I27").Value = "Mammals subtotal is inconsistent"

# Rows 29 and 31 (2009 and 2010, table 3): note work in progress cleaning the csvs.
$ws.Range("I29").Value = "Working on cleaning up"
$ws.Range("I31").Value = "Working on cleaning up"

# Scroll the window down (row 10 at the top) and leave the selection on I32,
# matching where the author ended up.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("I32").Select()
